$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2808.8235
$ws.Range("J64").Value = 3037.5
$ws.Range("L64").Value = 3037.5
$ws.Range("N64").Value = -3533.5

$ws.Range("H67").Value = 2808.8235
$ws.Range("J67").Value = 3037.5
$ws.Range("L67").Value = 3037.5
$ws.Range("N67").Value = -4753.5

$ws.Range("H92").Value = 734.82355
$ws.Range("I92").Value = 522.4545000000001
$ws.Range("J92").Value = 1124.1666
$ws.Range("K92").Value = 522.4545000000001
$ws.Range("L92").Value = 1124.1666
$ws.Range("M92").Value = 725.5454999999999
$ws.Range("N92").Value = -3620.1666

$ws.Range("H113").Value = 1980.7667
$ws.Range("I113").Value = 1835.5883
$ws.Range("K113").Value = 1835.5883
$ws.Range("M113").Value = 1418.4117

$ws.Range("H132").Value = 4927829.5
$ws.Range("I132").Value = 5292743
$ws.Range("J132").Value = 1495
$ws.Range("K132").Value = 15878229
$ws.Range("L132").Value = 4485
$ws.Range("M132").Value = -15875699
$ws.Range("N132").Value = -9545

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 22886
$ws.Range("J76").Value = 22886
$ws.Range("L76").Value = 22886
$ws.Range("N76").Value = -23562

$ws.Range("H79").Value = 22886
$ws.Range("J79").Value = 22886
$ws.Range("L79").Value = 22886
$ws.Range("N79").Value = -25226

$ws.Range("H92").Value = 40683.332
$ws.Range("J92").Value = 40683.332
$ws.Range("L92").Value = 40683.332
$ws.Range("N92").Value = -45675.332

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H106").Value = 41262.5
$ws.Range("J106").Value = 41262.5
$ws.Range("L106").Value = 41262.5
$ws.Range("N106").Value = -43786.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 23657
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 23657
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -24287

$ws.Range("H79").Value = 23657
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 23657
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -25841

$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 29500
$ws.Range("J74").Value = 29500
$ws.Range("L74").Value = 29500
$ws.Range("N74").Value = -31248

$ws.Range("H77").Value = 29500
$ws.Range("J77").Value = 29500
$ws.Range("L77").Value = 88500
$ws.Range("N77").Value = -97236

$ws.Range("H92").Value = 19927.715
$ws.Range("J92").Value = 19927.715
$ws.Range("L92").Value = 19927.715
$ws.Range("N92").Value = -24919.715

$ws.Range("H99").Value = 2220.7334
$ws.Range("I99").Value = 1812.3334
$ws.Range("K99").Value = 1812.3334
$ws.Range("M99").Value = -314.3334

$ws.Range("H107").Value = 433.52777
$ws.Range("I107").Value = 402.95456
$ws.Range("J107").Value = 481.57144
$ws.Range("K107").Value = 402.95456
$ws.Range("L107").Value = 481.57144
$ws.Range("M107").Value = 1517.04544
$ws.Range("N107").Value = -4321.57144

$ws.Range("H126").Value = 2220.7334
$ws.Range("I126").Value = 1812.3334
$ws.Range("K126").Value = 5437.0002
$ws.Range("M126").Value = -2967.0002

$ws.Range("H135").Value = 34780
$ws.Range("J135").Value = 34780
$ws.Range("L135").Value = 34780
$ws.Range("N135").Value = -44920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 22.625
$ws.Range("I2").Value = 27.363636
$ws.Range("J2").Value = 12.2
$ws.Range("K2").Value = 164.181816
$ws.Range("L2").Value = 73.19999999999999
$ws.Range("M2").Value = -51.181816
$ws.Range("N2").Value = -299.2

$ws.Range("H99").Value = 10219.917
$ws.Range("I99").Value = 1128
$ws.Range("J99").Value = 16714.143
$ws.Range("K99").Value = 3384
$ws.Range("L99").Value = 50142.429
$ws.Range("M99").Value = -1138
$ws.Range("N99").Value = -54634.429

$ws.Range("H118").Value = 3874.875
$ws.Range("I118").Value = 500
$ws.Range("J118").Value = 4999.8335
$ws.Range("K118").Value = 1500
$ws.Range("L118").Value = 14999.5005
$ws.Range("M118").Value = -257
$ws.Range("N118").Value = -17485.5005

$ws.Range("H131").Value = 3708113.8
$ws.Range("I131").Value = 27807.5
$ws.Range("J131").Value = 4274314.5
$ws.Range("K131").Value = 83422.5
$ws.Range("L131").Value = 12822943.5
$ws.Range("M131").Value = -78382.5
$ws.Range("N131").Value = -12833023.5

$ws.Range("H140").Value = 1562.3334
$ws.Range("I140").Value = 1244.9445
$ws.Range("J140").Value = 3466.6667
$ws.Range("K140").Value = 3734.8335
$ws.Range("L140").Value = 10400.0001
$ws.Range("M140").Value = 1445.1665
$ws.Range("N140").Value = -20760.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1915
$ws.Range("I126").Value = 1493.3334
$ws.Range("J126").Value = 3180
$ws.Range("K126").Value = 4480.0002
$ws.Range("L126").Value = 9540
$ws.Range("M126").Value = -2010.0002
$ws.Range("N126").Value = -14480

$ws.Range("H132").Value = 108270.42
$ws.Range("I132").Value = 202384
$ws.Range("J132").Value = 3699.7778
$ws.Range("K132").Value = 607152
$ws.Range("L132").Value = 11099.3334
$ws.Range("M132").Value = -604622
$ws.Range("N132").Value = -16159.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1998
$ws.Range("I82").Value = 2060
$ws.Range("K82").Value = 2060
$ws.Range("M82").Value = -1699

$ws.Range("H85").Value = 1998
$ws.Range("I85").Value = 2060
$ws.Range("K85").Value = 2060
$ws.Range("M85").Value = -812

$ws.Range("H106").Value = 20000
$ws.Range("J106").Value = 20000
$ws.Range("L106").Value = 20000
$ws.Range("N106").Value = -22524

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 78476.336
$ws.Range("J138").Value = 78476.336
$ws.Range("L138").Value = 78476.336
$ws.Range("N138").Value = -88756.336
